$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 100
$ws.Range("I9").Value = 100
$ws.Range("K9").Value = 100
$ws.Range("M9").Value = 69

$ws.Range("H62").Value = 43120.883
$ws.Range("I62").Value = 59519.668
$ws.Range("J62").Value = 6223.625
$ws.Range("K62").Value = 59519.668
$ws.Range("L62").Value = 6223.625
$ws.Range("M62").Value = -58895.668
$ws.Range("N62").Value = -7471.625

$ws.Range("H65").Value = 43120.883
$ws.Range("I65").Value = 59519.668
$ws.Range("J65").Value = 6223.625
$ws.Range("K65").Value = 297598.34
$ws.Range("L65").Value = 31118.125
$ws.Range("M65").Value = -294478.34
$ws.Range("N65").Value = -37358.125

$ws.Range("H129").Value = 1001.0351
$ws.Range("J129").Value = 1135.1063
$ws.Range("L129").Value = 3405.3189
$ws.Range("N129").Value = -13405.3189

$ws.Range("H137").Value = 700.6863
$ws.Range("I137").Value = 647.4667
$ws.Range("J137").Value = 1099.8334
$ws.Range("K137").Value = 1942.4001
$ws.Range("L137").Value = 3299.5002
$ws.Range("M137").Value = 607.5999000000002
$ws.Range("N137").Value = -8399.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 386.66666
$ws.Range("I5").Value = 80
$ws.Range("K5").Value = 80
$ws.Range("M5").Value = 32

$ws.Range("H14").Value = 100
$ws.Range("I14").Value = 100
$ws.Range("K14").Value = 100
$ws.Range("M14").Value = 75

$ws.Range("H32").Value = 6154.066
$ws.Range("I32").Value = 5955.281
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 5955.281
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -5668.281
$ws.Range("N32").Value = -15574

$ws.Range("H74").Value = 1000.13336
$ws.Range("I74").Value = 983.3939
$ws.Range("J74").Value = 1046.1666
$ws.Range("K74").Value = 983.3939
$ws.Range("L74").Value = 1046.1666
$ws.Range("M74").Value = -109.3939
$ws.Range("N74").Value = -2794.1666

$ws.Range("H77").Value = 1000.13336
$ws.Range("I77").Value = 983.3939
$ws.Range("J77").Value = 1046.1666
$ws.Range("K77").Value = 4916.9695
$ws.Range("L77").Value = 5230.833000000001
$ws.Range("M77").Value = -548.9695000000002
$ws.Range("N77").Value = -13966.833

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 386.66666
$ws.Range("I4").Value = 80
$ws.Range("K4").Value = 80
$ws.Range("M4").Value = 35

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2801240.8
$ws.Range("I2").Value = 5882490
$ws.Range("J2").Value = 105.36364
$ws.Range("K2").Value = 35294940
$ws.Range("L2").Value = 632.18184
$ws.Range("M2").Value = -35294827
$ws.Range("N2").Value = -858.18184

$ws.Range("H23").Value = 117.882355
$ws.Range("I23").Value = 78.2
$ws.Range("J23").Value = 134.41667
$ws.Range("K23").Value = 234.6
$ws.Range("L23").Value = 403.25001
$ws.Range("M23").Value = 0.3999999999999773
$ws.Range("N23").Value = -873.25001

$ws.Range("H40").Value = 193
$ws.Range("I40").Value = 190
$ws.Range("J40").Value = 194.5
$ws.Range("K40").Value = 760
$ws.Range("L40").Value = 778
$ws.Range("M40").Value = -691
$ws.Range("N40").Value = -916

$ws.Range("H46").Value = 4955.8125
$ws.Range("I46").Value = 829.3
$ws.Range("J46").Value = 11833.333
$ws.Range("K46").Value = 2487.9
$ws.Range("L46").Value = 35499.999
$ws.Range("M46").Value = -2396.9
$ws.Range("N46").Value = -35681.999

$ws.Range("H51").Value = 3180.6
$ws.Range("I51").Value = 2001.5
$ws.Range("J51").Value = 3966.6667
$ws.Range("K51").Value = 6004.5
$ws.Range("L51").Value = 11900.0001
$ws.Range("M51").Value = -5544.5
$ws.Range("N51").Value = -12820.0001

$ws.Range("H57").Value = 3200
$ws.Range("I57").Value = 2900
$ws.Range("J57").Value = 3500
$ws.Range("K57").Value = 8700
$ws.Range("L57").Value = 10500
$ws.Range("M57").Value = -8141
$ws.Range("N57").Value = -11618

$ws.Range("H58").Value = 2087.2727
$ws.Range("I58").Value = 650
$ws.Range("J58").Value = 2406.6667
$ws.Range("K58").Value = 1950
$ws.Range("L58").Value = 7220.000100000001
$ws.Range("M58").Value = -1822
$ws.Range("N58").Value = -7476.000100000001

$ws.Range("H69").Value = 312
$ws.Range("I69").Value = 312
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 936
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -125
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 312
$ws.Range("I72").Value = 312
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 2808
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = 1248
$ws.Range("N72").ClearContents()

$ws.Range("H86").Value = 600
$ws.Range("I86").Value = 500
$ws.Range("J86").Value = 800
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 2400
$ws.Range("M86").Value = -314
$ws.Range("N86").Value = -4772

$ws.Range("H87").Value = 9421.777
$ws.Range("I87").Value = 5049.3335
$ws.Range("J87").Value = 18166.666
$ws.Range("K87").Value = 15148.0005
$ws.Range("L87").Value = 54499.99800000001
$ws.Range("M87").Value = -13900.0005
$ws.Range("N87").Value = -56995.99800000001

$ws.Range("H89").Value = 600
$ws.Range("I89").Value = 500
$ws.Range("J89").Value = 800
$ws.Range("K89").Value = 4500
$ws.Range("L89").Value = 7200
$ws.Range("M89").Value = 1428
$ws.Range("N89").Value = -19056

$ws.Range("H90").Value = 9421.777
$ws.Range("I90").Value = 5049.3335
$ws.Range("J90").Value = 18166.666
$ws.Range("K90").Value = 45444.0015
$ws.Range("L90").Value = 163499.994
$ws.Range("M90").Value = -39204.0015
$ws.Range("N90").Value = -175979.994

$ws.Range("H137").Value = 48938.773
$ws.Range("J137").Value = 129675
$ws.Range("L137").Value = 389025
$ws.Range("N137").Value = -399225

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 33334542
$ws.Range("I46").Value = 40001260
$ws.Range("J46").Value = 950
$ws.Range("K46").Value = 40001260
$ws.Range("L46").Value = 950
$ws.Range("M46").Value = -40001072
$ws.Range("N46").Value = -1326

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 33500
$ws.Range("J114").Value = 33500
$ws.Range("L114").Value = 33500
$ws.Range("N114").Value = -42178
